# Adds a "Tools" sheet (GitHub git-leak tooling links) right before the
# "Udemy" sheet, and an "Article" sheet (FreeCodeCamp SQL-injection article)
# right after "Udemy" — per commit "FreeCodeCamp article about SQLi added."

$wb = $excel.ActiveWorkbook
$udemy = $wb.Worksheets.Item("Udemy")

# ---------------------------------------------------------------------
# New "Tools" worksheet, inserted immediately before "Udemy"
# ---------------------------------------------------------------------
$tools = $wb.Worksheets.Add($udemy)
$tools.Name = "Tools"

# Fill column C (urls) before column B (repeated label) so the shared
# string table is populated in the same order as the source workbook.
$tools.Range("C2").Value = "https://github.com/dnoiz1/git-money"
$tools.Range("C4").Value = "https://github.com/evilpacket/DVCS-Pillage"
$tools.Range("C6").Value = "https://github.com/internetwache/GitTools"

$tools.Range("B2").Value = "GIT leak Tools"
$tools.Range("B4").Value = "GIT leak Tools"
$tools.Range("B6").Value = "GIT leak Tools"

$tools.Columns.Item(2).ColumnWidth = 70.259
$tools.Columns.Item(3).ColumnWidth = 61.259

[void]$tools.Range("B8").Select()

# ---------------------------------------------------------------------
# New "Article" worksheet, inserted immediately after "Udemy"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$article = $wb.Worksheets.Add($null, $lastSheet)
$article.Name = "Article"

$article.Range("C2").Value = "https://www.freecodecamp.org/news/what-is-sql-injection-how-to-prevent-it/"
$article.Range("B2").Value = "FreeCodeCamp - SQL Injection "

$article.Columns.Item(2).ColumnWidth = 70.09
$article.Columns.Item(3).ColumnWidth = 96.754

[void]$article.Range("B2").Select()
[void]$article.Activate()
